# Pin Assignments workbook update:
#  - Add a new "P3" connector block (rows 73-96) describing the FPGA-based
#    velocity-measurement inputs for the new encoders, following the same
#    layout/formatting already used for the P4 (rows 1-27) and P5
#    (rows 28-51) connector blocks above.
#  - Row 72 (formerly a blank filler row) becomes the thick-bottom-border
#    closing row of the preceding block.
#  - Two new cell formats are introduced (right-aligned versions of the
#    existing "header" fills) for the new P3 header row's Pin column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Close out the previous block: row 72 gets the thick-bottom style (style
#    used by rows 27/51, i.e. fillId=0/borderId=medium) and the matching
#    slightly-taller row height used by the other section-ending rows.
# ---------------------------------------------------------------------------
$ws.Range("A1:F1").Copy()
$ws.Range("A72:F72").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(72).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 2) New row 73 ("P3" header / Ground pin) - same formatting as row 28 (the
#    "P5" header row), with the Pin (C) column additionally right-aligned.
# ---------------------------------------------------------------------------
$ws.Range("A28:F28").Copy()
$ws.Range("A73:F73").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C73").HorizontalAlignment = -4152  # xlRight

# ---------------------------------------------------------------------------
# 3) New row 74 (5V pin) - same formatting as row 3, with the Pin (C) column
#    additionally right-aligned.
# ---------------------------------------------------------------------------
$ws.Range("A3:F3").Copy()
$ws.Range("A74:F74").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C74").HorizontalAlignment = -4152  # xlRight

# ---------------------------------------------------------------------------
# 4) New rows 75-95: plain data rows, same formatting as row 4 (style 2
#    throughout).
# ---------------------------------------------------------------------------
$ws.Range("A4:F4").Copy()
$ws.Range("A75:F95").PasteSpecial(-4122)  # xlPasteFormats

# A couple of individual cells within that range deviate from the plain
# style-2 formatting, matching the existing convention used elsewhere in the
# sheet (e.g. D33) of highlighting specific pins with the header fill:
$ws.Range("D28").Copy()
$ws.Range("D78").PasteSpecial(-4122)
$ws.Range("A28").Copy()
$ws.Range("E86").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5) New row 96: closing row of the P3 block, same formatting as row 27/51
#    (thick bottom border) plus the matching row height.
# ---------------------------------------------------------------------------
$ws.Range("A27:F27").Copy()
$ws.Range("A96:F96").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(96).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 6) Fill in the values. Row 73/74 "E" column values are written in this
#    particular order so the new shared strings are appended in the same
#    order as in the target workbook (P3, then "5V+ rail...", then
#    "Ground Rail...").
# ---------------------------------------------------------------------------
$ws.Range("A73").Value2 = "P3"
$ws.Range("B73").Value2 = 50
$ws.Range("C73").Value2 = "x"
$ws.Range("D73").Value2 = "D GND"
$ws.Range("F73").Value2 = "FPGA"

$ws.Range("B74").Value2 = 48
$ws.Range("C74").Value2 = "x"
$ws.Range("D74").Value2 = "5V"
$ws.Range("F74").Value2 = "FPGA"

$ws.Range("E74").Value2 = "5V+ rail for encoders"
$ws.Range("E73").Value2 = "Ground Rail for encoders"

$ws.Range("B75").Value2 = 45
$ws.Range("C75").Value2 = 9
$ws.Range("D75").Value2 = "DIO6"
$ws.Range("F75").Value2 = "FPGA"

$ws.Range("B76").Value2 = 9
$ws.Range("C76").Value2 = 7
$ws.Range("D76").Value2 = "DIO3"
$ws.Range("F76").Value2 = "FPGA"

$ws.Range("B77").Value2 = 11
$ws.Range("C77").Value2 = 3
$ws.Range("D77").Value2 = "DIO9"
$ws.Range("F77").Value2 = "FPGA"

$ws.Range("B78").Value2 = 13
$ws.Range("C78").Value2 = 3
$ws.Range("D78").Value2 = "DIO0"
$ws.Range("F78").Value2 = "FPGA"

$ws.Range("B79").Value2 = 15
$ws.Range("C79").Value2 = 3
$ws.Range("D79").Value2 = "DIO1"
$ws.Range("F79").Value2 = "FPGA"

$ws.Range("B80").Value2 = 17
$ws.Range("C80").Value2 = 3
$ws.Range("D80").Value2 = "DIO2"
$ws.Range("F80").Value2 = "FPGA"

$ws.Range("B81").Value2 = 19
$ws.Range("C81").Value2 = 3
$ws.Range("D81").Value2 = "DIO3"
$ws.Range("F81").Value2 = "FPGA"

$ws.Range("B82").Value2 = 21
$ws.Range("C82").Value2 = 3
$ws.Range("D82").Value2 = "DIO4"
$ws.Range("F82").Value2 = "FPGA"

$ws.Range("B83").Value2 = 23
$ws.Range("C83").Value2 = 3
$ws.Range("D83").Value2 = "DIO5"
$ws.Range("F83").Value2 = "FPGA"

$ws.Range("B84").Value2 = 25
$ws.Range("C84").Value2 = 3
$ws.Range("D84").Value2 = "DIO6"
$ws.Range("F84").Value2 = "FPGA"

$ws.Range("B85").Value2 = 27
$ws.Range("C85").Value2 = 3
$ws.Range("D85").Value2 = "DIO7"
$ws.Range("F85").Value2 = "FPGA"

$ws.Range("B86").Value2 = 29
$ws.Range("C86").Value2 = 3
$ws.Range("D86").Value2 = "DIO8"
$ws.Range("F86").Value2 = "FPGA"

$ws.Range("B87").Value2 = 31
$ws.Range("C87").Value2 = 4
$ws.Range("D87").Value2 = "DIO9"
$ws.Range("F87").Value2 = "FPGA"

$ws.Range("B88").Value2 = 33
$ws.Range("C88").Value2 = 4
$ws.Range("D88").Value2 = "DIO0"
$ws.Range("F88").Value2 = "FPGA"

$ws.Range("B89").Value2 = 35
$ws.Range("C89").Value2 = 4
$ws.Range("D89").Value2 = "DIO1"
$ws.Range("F89").Value2 = "FPGA"

$ws.Range("B90").Value2 = 37
$ws.Range("C90").Value2 = 4
$ws.Range("D90").Value2 = "DIO2"
$ws.Range("F90").Value2 = "FPGA"

$ws.Range("B91").Value2 = 39
$ws.Range("C91").Value2 = 4
$ws.Range("D91").Value2 = "DIO3"
$ws.Range("F91").Value2 = "FPGA"

$ws.Range("B92").Value2 = 41
$ws.Range("C92").Value2 = 4
$ws.Range("D92").Value2 = "DIO4"
$ws.Range("F92").Value2 = "FPGA"

$ws.Range("B93").Value2 = 43
$ws.Range("C93").Value2 = 4
$ws.Range("D93").Value2 = "DIO5"
$ws.Range("F93").Value2 = "FPGA"

$ws.Range("B94").Value2 = 45
$ws.Range("C94").Value2 = 4
$ws.Range("D94").Value2 = "DIO6"
$ws.Range("F94").Value2 = "FPGA"

$ws.Range("B95").Value2 = 47
$ws.Range("C95").Value2 = 4
$ws.Range("D95").Value2 = "DIO7"
$ws.Range("F95").Value2 = "FPGA"

$ws.Range("B96").Value2 = 49
$ws.Range("C96").Value2 = 4
$ws.Range("D96").Value2 = "DIO8"
$ws.Range("F96").Value2 = "FPGA"

# ---------------------------------------------------------------------------
# 7) Update the active selection to match where the author ended up working.
# ---------------------------------------------------------------------------
$ws.Range("C75").Select()
